# Lab 8 hint edit: replace Hint 2 text, and move the lastRenderedPageBreak
# rendering marker to reflect the new pagination (removed from the
# "Your job is to modify your Die class..." run, added before
# "Integer.toString" later in the document).

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like $needle) {
            return $i
        }
    }
    return -1
}

# 1) "Hint 2: ..." paragraph - collapse the multi-run explanation + code
#    fragment into a single plain-text hint about order of operations.
$hintIdx = Find-ParagraphIndex $d "*Hint 2:*drawSpot*"
if ($hintIdx -lt 0) { throw "Could not locate the Hint 2 paragraph" }
$hintRange = $d.Paragraphs.Item($hintIdx).Range
$hintRange.Find.Execute(
    "You must have the division part of the scale at the very end of the math or you will lose values In the Int conversion for drawSpot().",
    $true, $false, $false, $false, $false, $true, 1, $false,
    'When using a multiplication factor, be careful about what order you do the math and your brackets since you can easily run Into unit conversion loss which can throw off the stretch.', 2) | Out-Null

# 2) Drop the stale <w:lastRenderedPageBreak/> that was sitting on the
#    "Your job is to modify your Die class..." paragraph.
$dieIdx = Find-ParagraphIndex $d "*Your job is to modify your*class to make those things happen*"
if ($dieIdx -lt 0) { throw "Could not locate the 'Your job is to modify your' paragraph" }
$d.Paragraphs.Item($dieIdx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0D91E40C" w14:textId="099ECE35" w:rsidR="00D57D08" w:rsidRDefault="00B02A6B" w:rsidP="00B02A6B"><w:r><w:t xml:space="preserve">Your job is to modify your </w:t></w:r><w:r w:rsidRPr="00B02A6B"><w:rPr><w:rFonts w:ascii="CourierNewPSMT" w:hAnsi="CourierNewPSMT"/><w:szCs w:val="26"/></w:rPr><w:t>Die</w:t></w:r><w:r><w:t xml:space="preserve"> class to make those things happen.</w:t></w:r></w:p>')

# 3) Re-add the marker where the page now actually breaks, just before
#    the "Integer.toString" run later in the document.
$introIdx = Find-ParagraphIndex $d "*must be converted to a String*Integer.toString*may help you*"
if ($introIdx -lt 0) { throw "Could not locate the 'Integer.toString' paragraph" }
$d.Paragraphs.Item($introIdx).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1B90BB66" w14:textId="51F5434D" w:rsidR="007D4D9E" w:rsidRDefault="007D4D9E" w:rsidP="007D4D9E"><w:pPr><w:widowControl w:val="0"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="360"/></w:tabs><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00F42233"><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F42233"><w:rPr><w:rFonts w:ascii="CourierNewPSMT" w:hAnsi="CourierNewPSMT"/><w:szCs w:val="26"/></w:rPr><w:t>currentValue</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F42233"><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> should be printed as </w:t></w:r><w:r w:rsidRPr="00F42233"><w:rPr><w:rFonts w:ascii="TimesNewRomanPS-ItalicMT" w:hAnsi="TimesNewRomanPS-ItalicMT"/><w:i/><w:szCs w:val="32"/></w:rPr><w:t>black</w:t></w:r><w:r w:rsidRPr="00F42233"><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> text somewhere inside the die.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">  To do so, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F42233"><w:rPr><w:rFonts w:ascii="CourierNewPSMT" w:hAnsi="CourierNewPSMT"/><w:szCs w:val="26"/></w:rPr><w:t>currentValue</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> must be converted to a String.</w:t></w:r><w:r w:rsidR="004C079B"><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004C079B" w:rsidRPr="004C079B"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>Integer.toString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004C079B" w:rsidRPr="004C079B"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:szCs w:val="32"/></w:rPr><w:t>(int)</w:t></w:r><w:r w:rsidR="004C079B"><w:rPr><w:rFonts w:ascii="TimesNewRomanPSMT" w:hAnsi="TimesNewRomanPSMT"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> may help you. </w:t></w:r></w:p>')

Write-Host "Done."
